# Business Request - Mail from Steven.docx
# Change the budget year mentioned in the document from 2021 to 2023.
#
# The author made this edit the way a human does in Word: click right
# before the final "1" in "2021", select just that one character, and
# type "3" over it. That kind of in-place, single-character retype
# leaves the paragraph split into three runs (the text before the
# retyped character, the retyped character itself, and the text after
# it) even though all three runs end up with identical run properties.
#
# Reproduce that: locate "2021", narrow down to its last character
# ("1"), nudge its direct character formatting away from the
# surrounding text so it is recorded as its own run, overwrite its
# text with "3", then restore the formatting so the final run
# properties match the rest of the paragraph.

$d = $word.ActiveDocument

# Find the "2021" occurrence; after Execute returns, $yearRange is
# collapsed onto the matched text ("2021").
$yearRange = $d.Content
$yearRange.Find.Execute("2021", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

# Narrow down to just the last character of the match (the "1").
$digitStart = $yearRange.End - 1
$digitEnd = $yearRange.End
$digitRange = $d.Range($digitStart, $digitEnd)

# Make this character its own run (temporarily bold it), retype it,
# then drop the bold again so the formatting ends up unchanged.
$digitRange.Bold = 1
$digitRange.Text = "3"
$digitRange.Bold = 0
